{"js": "// 1) Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\"\nconst dateResults = context.document.body.search(\"September 19, 2025\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"September 21, 2025\", Word.InsertLocation.replace);\n}\n\n// 2) Split the mailing-address line \"1682 East Avenue, Fairview CA 94541\" into two\n//    paragraphs: \"1682 East Avenue\" and a new paragraph \"Fairview, CA 94541\".\n//    (Only the standalone address paragraph near the top of the letter is touched,\n//    not the identical text that also appears inside the PROPERTY ADDRESS table cell.)\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet addressParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"1682 East Avenue, Fairview CA 94541\") {\n    addressParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (addressParagraph) {\n  // Insert the new second line right after the address paragraph; it inherits the\n  // same paragraph/run formatting (Arial 11pt, autoSpaceDE/autoSpaceDN off).\n  addressParagraph.insertParagraph(\"Fairview, CA 94541\", Word.InsertLocation.after);\n  // Trim the original paragraph's text down to just the street address.\n  addressParagraph.getRange().insertText(\"1682 East Avenue\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n\n// 3) Remove the empty \"No Spacing\" paragraph that sits right after the\n//    \"... Board of Directors\" signature line.\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (paragraphs2.items[i].text === \"Shady Hollow Owners Association Board of Directors\") {\n    const next = paragraphs2.items[i + 1];\n    if (next) {\n      next.load(\"text\");\n      await context.sync();\n      if (next.text === \"\") {\n        next.delete();\n      }\n    }\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\"\n$dateRange = $d.Content\n$dateFind = $dateRange.Find\n$dateFind.Text = \"September 19, 2025\"\n$dateFound = $dateFind.Execute()\nif ($dateFound) {\n    $dateRange.Text = \"September 21, 2025\"\n}\n\n# 2) Split the mailing-address line \"1682 East Avenue, Fairview CA 94541\" into two\n#    paragraphs: \"1682 East Avenue\" and a new paragraph \"Fairview, CA 94541\".\n#    Only the standalone address paragraph near the top of the letter is touched,\n#    not the identical text that also appears inside the PROPERTY ADDRESS table cell.\n$addressIndex = 0\n$i = 1\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"1682 East Avenue, Fairview CA 94541\") {\n        $addressIndex = $i\n        break\n    }\n    $i = $i + 1\n}\n\nif ($addressIndex -gt 0) {\n    $addrPara = $d.Paragraphs($addressIndex)\n    $addrPara.Range.InsertParagraphAfter()\n    $d.Paragraphs($addressIndex).Range.Text = \"1682 East Avenue\"\n    $d.Paragraphs($addressIndex + 1).Range.Text = \"Fairview, CA 94541\"\n}\n\n# 3) Remove the empty \"No Spacing\" paragraph that sits right after the\n#    \"... Board of Directors\" signature line.\n$boardIndex = 0\n$i = 1\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"Shady Hollow Owners Association Board of Directors\") {\n        $boardIndex = $i\n        break\n    }\n    $i = $i + 1\n}\n\nif ($boardIndex -gt 0) {\n    $nextPara = $d.Paragraphs($boardIndex + 1)\n    $nextText = $nextPara.Range.Text.TrimEnd([char]13, [char]7)\n    if ($nextText -eq \"\") {\n        $nextPara.Range.Delete()\n    }\n}\n"}
